$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "II Scenario/Flusso di eventi di ERRORE" row (error scenario header)
# and the "Sistema:" row directly below it (which holds the error-message text),
# by scanning the table instead of relying on a hard-coded row index.
$headerRow = $null
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $cellText = $t.Rows.Item($r).Cells.Item(1).Range.Text
    if ($cellText -like "*II Scenario/Flusso di eventi di ERRORE*sistema non riesce a comunicare col database*") {
        $headerRow = $r
        break
    }
}

if ($headerRow -eq $null) {
    throw "Could not locate the 'II Scenario/Flusso di eventi di ERRORE' row"
}

# --- Hunk 1: fix the scenario header wording ---
$headerCell = $t.Rows.Item($headerRow).Cells.Item(1)
$headerCell.Range.Find.Execute(
    "sistema non riesce a comunicare col database", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "il sistema non riesce ad effettuare il salvataggio dei dati", 2) | Out-Null

# --- Hunk 2: fix the error-message description in the following "Sistema:" row ---
$msgRow = $headerRow + 1
$msgCell = $t.Rows.Item($msgRow).Cells.Item(3)
$msgCell.Range.Find.Execute(
    "Visualizza un messaggio di errore all'utente. Il messaggio segnala non è stato possibile effettuare la comunicazione di sistema.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Visualizza un messaggio di errore al client. Il messaggio segnala che non è stato possibile salvare i dati e invita a riprovare più tardi.", 2) | Out-Null
